# Refresh cached Market Board price/profit figures in the Leviathan Profits
# workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve sheets), as produced by the
# scheduled price-update runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H51").Value = 5559484
$ws.Range("J51").Value = 16670467
$ws.Range("L51").Value = 16670467
$ws.Range("N51").Value = -16671435

$ws.Range("H98").Value = 1896.0769
$ws.Range("I98").Value = 1434.5714
$ws.Range("J98").Value = 3834.4
$ws.Range("K98").Value = 1434.5714
$ws.Range("L98").Value = 3834.4
$ws.Range("M98").Value = 63.42859999999996
$ws.Range("N98").Value = -6830.4

$ws.Range("H104").Value = 245.33333
$ws.Range("I104").Value = 245.33333
$ws.Range("K104").Value = 735.99999
$ws.Range("M104").Value = 1011.00001

$ws.Range("H122").Value = 1896.0769
$ws.Range("I122").Value = 1434.5714
$ws.Range("J122").Value = 3834.4
$ws.Range("K122").Value = 4303.7142
$ws.Range("L122").Value = 11503.2
$ws.Range("M122").Value = -1853.7142
$ws.Range("N122").Value = -16403.2

$ws.Range("H132").Value = 1803.1389
$ws.Range("I132").Value = 1621
$ws.Range("K132").Value = 4863
$ws.Range("M132").Value = -2333

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 2164.739
$ws.Range("I2").Value = 1804.4762
$ws.Range("J2").Value = 5947.5
$ws.Range("K2").Value = 1804.4762
$ws.Range("L2").Value = 5947.5
$ws.Range("M2").Value = -1691.4762
$ws.Range("N2").Value = -6173.5

$ws.Range("H45").Value = 5299.485
$ws.Range("I45").Value = 6416.381
$ws.Range("K45").Value = 6416.381
$ws.Range("M45").Value = -6039.381

$ws.Range("H74").Value = 2044.35
$ws.Range("I74").Value = 1924.2778
$ws.Range("K74").Value = 1924.2778
$ws.Range("M74").Value = -1050.2778

$ws.Range("H77").Value = 2044.35
$ws.Range("I77").Value = 1924.2778
$ws.Range("K77").Value = 9621.389000000001
$ws.Range("M77").Value = -5253.389000000001

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H116").Value = 2164.739
$ws.Range("I116").Value = 1804.4762
$ws.Range("J116").Value = 5947.5
$ws.Range("K116").Value = 1804.4762
$ws.Range("L116").Value = 5947.5
$ws.Range("M116").Value = 489.5237999999999
$ws.Range("N116").Value = -10535.5

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H3").Value = 2164.739
$ws.Range("I3").Value = 1804.4762
$ws.Range("J3").Value = 5947.5
$ws.Range("K3").Value = 1804.4762
$ws.Range("L3").Value = 5947.5
$ws.Range("M3").Value = -1690.4762
$ws.Range("N3").Value = -6175.5

$ws.Range("H38").Value = 40000
$ws.Range("J38").Value = 40000
$ws.Range("L38").Value = 40000
$ws.Range("N38").Value = -40832

$ws.Range("H42").Value = 163882
$ws.Range("J42").Value = 163882
$ws.Range("L42").Value = 163882
$ws.Range("N42").Value = -164538

$ws.Range("H105").Value = 826.03845
$ws.Range("I105").Value = 800.5454999999999
$ws.Range("K105").Value = 800.5454999999999
$ws.Range("M105").Value = 946.4545000000001

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H3").Value = 15555794
$ws.Range("I3").Value = 17500254
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 17500254
$ws.Range("L3").Value = 120
$ws.Range("M3").Value = -17500141
$ws.Range("N3").Value = -346

$ws.Range("H86").Value = 4141.4287
$ws.Range("I86").Value = 3800.4
$ws.Range("K86").Value = 3800.4
$ws.Range("M86").Value = -2677.4

$ws.Range("H89").Value = 4141.4287
$ws.Range("I89").Value = 3800.4
$ws.Range("K89").Value = 19002
$ws.Range("M89").Value = -13386

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H2").Value = 463.35715
$ws.Range("I2").Value = 886
$ws.Range("J2").Value = 40.714287
$ws.Range("K2").Value = 5316
$ws.Range("L2").Value = 244.285722
$ws.Range("M2").Value = -5203
$ws.Range("N2").Value = -470.285722

$ws.Range("H34").Value = 180.53847
$ws.Range("J34").Value = 1000
$ws.Range("L34").Value = 3000
$ws.Range("N34").Value = -3168

$ws.Range("H39").Value = 8294.474
$ws.Range("J39").Value = 8294.474
$ws.Range("L39").Value = 24883.422
$ws.Range("N39").Value = -25471.422

$ws.Range("H55").Value = 6581377.5
$ws.Range("I55").Value = 911.6
$ws.Range("J55").Value = 31258124
$ws.Range("K55").Value = 2734.8
$ws.Range("L55").Value = 93774372
$ws.Range("M55").Value = -2557.8
$ws.Range("N55").Value = -93774726

$ws.Range("H113").Value = 743
$ws.Range("I113").Value = 724.5
$ws.Range("K113").Value = 2173.5
$ws.Range("M113").Value = -3.5

$ws.Range("H122").Value = 382.25
$ws.Range("I122").Value = 454
$ws.Range("K122").Value = 4086
$ws.Range("M122").Value = -1636

$ws.Range("H129").Value = 1844.0416
$ws.Range("I129").Value = 509
$ws.Range("J129").Value = 3713.1
$ws.Range("K129").Value = 1527
$ws.Range("L129").Value = 11139.3
$ws.Range("M129").Value = 3473
$ws.Range("N129").Value = -21139.3

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H80").Value = 5199.5
$ws.Range("I80").Value = 3399.8
$ws.Range("K80").Value = 3399.8
$ws.Range("M80").Value = -2401.8

$ws.Range("H83").Value = 5199.5
$ws.Range("I83").Value = 3399.8
$ws.Range("K83").Value = 16999
$ws.Range("M83").Value = -12007

$ws.Range("H126").Value = 2883.8948
$ws.Range("I126").Value = 2532.9333
$ws.Range("K126").Value = 7598.7999
$ws.Range("M126").Value = -5128.7999

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 1275.3846

$ws.Range("H27").Value = 1275.3846

$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H81").Value = 14577
$ws.Range("I81").Value = 962
$ws.Range("K81").Value = 1924
$ws.Range("M81").Value = -863

$ws.Range("H84").Value = 14577
$ws.Range("I84").Value = 962
$ws.Range("K84").Value = 9620
$ws.Range("M84").Value = -4316

$ws.Range("H131").Value = 92160
$ws.Range("J131").Value = 92160
$ws.Range("L131").Value = 92160
$ws.Range("N131").Value = -102240
